$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: near the end of the document, drop the existing "_GoBack" bookmark
# that currently sits between "...betere optie." and the trailing space, and
# merge the trailing space into the preceding run so the sentence ends up as
# a single run: ", ... een betere optie. " (xml:space="preserve").
# This must run BEFORE we add a new "_GoBack" bookmark elsewhere, since
# bookmark names are unique in the document.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBackPos = $goBack.Range.Start
$goBack.Delete()

$contentEnd = $d.Content.End
$trailingSpace = $d.Range($goBackPos, ($contentEnd - 1))
$trailingSpace.Delete()

# Rewrite the final character of the remaining run together with a restored
# trailing space so both pieces of text end up inside the same run/xml:t.
$mergePoint = $d.Range(($goBackPos - 1), $goBackPos)
$mergePoint.Text = ". "

# ---------------------------------------------------------------------------
# Part 2: title paragraph ("Conclusie: native apps versus HTML5-apps") gets
# bumped to 14pt (sz/szCs = 28 half-points) on the paragraph mark and on
# every run.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1).Range
$titlePara.Font.Size = 14
$titlePara.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# Part 3: add a fresh, collapsed "_GoBack" bookmark right after the last run
# of the title paragraph (before the paragraph mark). Inserting it directly
# at that paragraph-final offset is ambiguous in this runtime, so a
# temporary placeholder character is used to get a stable anchor, the
# bookmark is created there, and the placeholder is removed again.
# ---------------------------------------------------------------------------
$titleEnd = $titlePara.End - 1
$placeholder = $d.Range($titleEnd, $titleEnd)
$placeholder.InsertAfter("X")

$bkRange = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bkRange)

$cleanup = $d.Range($titleEnd, ($titleEnd + 1))
$cleanup.Delete()
